$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: locate the last paragraph ("Preguntar si en la teoria de
# regresion logistica es necesario meter lo del logit.") and split its
# single run into three runs: "...meter lo del ", "logit", "."
# (the surrounding spell-check proofErr marks are cosmetic artifacts
# of Word's spell checker and are not reachable through the COM
# object model, so we focus on reproducing the run/text structure).
# ------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$p1 = $d.Paragraphs.Item($lastParaIndex)
$p1Start = $p1.Range.Start
$p1Text = $p1.Range.Text

$word1 = "logit"
$idxWord1 = $p1Text.IndexOf($word1)

# Range covering "logit"
$rLogit = $d.Range($p1Start + $idxWord1, $p1Start + $idxWord1 + $word1.Length)
$rLogit.Bold = 1
$rLogit.Bold = 0

# Range covering the trailing "." right after "logit"
$rDot = $d.Range($p1Start + $idxWord1 + $word1.Length, $p1Start + $idxWord1 + $word1.Length + 1)
$rDot.Bold = 1
$rDot.Bold = 0

# ------------------------------------------------------------------
# Step 2: append a new list paragraph (same style/numbering as the
# paragraph above, inherited automatically by InsertParagraphAfter)
# with the "Mirar tesis de Hanen Borchani..." text, split into runs
# around "Hanen" and "Borchani".
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p1.Range.InsertParagraphAfter() | Out-Null

$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$newText = "Mirar tesis de Hanen Borchani para definir mejor el problema de clasificación supervisada. -> Preguntar si es necesario definirlo mejor puesto que al nombrar la notación no sabemos si hay que mencionar lo que es cada cosa al establecer el problema de clasificación supervisada."
$p2.Range.Text = $newText

$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p2Start = $p2.Range.Start
$p2Text = $p2.Range.Text

$word2 = "Hanen"
$idxWord2 = $p2Text.IndexOf($word2)
$rHanen = $d.Range($p2Start + $idxWord2, $p2Start + $idxWord2 + $word2.Length)
$rHanen.Bold = 1
$rHanen.Bold = 0

# space right after "Hanen"
$rSpace = $d.Range($p2Start + $idxWord2 + $word2.Length, $p2Start + $idxWord2 + $word2.Length + 1)
$rSpace.Bold = 1
$rSpace.Bold = 0

$word3 = "Borchani"
$idxWord3 = $p2Text.IndexOf($word3)
$rBorchani = $d.Range($p2Start + $idxWord3, $p2Start + $idxWord3 + $word3.Length)
$rBorchani.Bold = 1
$rBorchani.Bold = 0

# remaining tail starting with the space right after "Borchani"
$rTail = $d.Range($p2Start + $idxWord3 + $word3.Length, $p2.Range.End - 1)
$rTail.Bold = 1
$rTail.Bold = 0

# ------------------------------------------------------------------
# Step 3: move the "_GoBack" bookmark so that it again sits right
# after the very last run, before the final paragraph mark (this is
# where Word leaves it after the most recent edit). The runtime has a
# quirk where Bookmarks.Add placed exactly at Content.End-1 resolves
# incorrectly, so a temporary trailing character is used to dodge
# that edge case and is removed again afterwards.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}

$sentinelPos = $d.Content.End - 1
$sentinelRange = $d.Range($sentinelPos, $sentinelPos)
$sentinelRange.InsertAfter("#")

$bmPos = $d.Content.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanupRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$cleanupRange.Delete()
